$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the contents of columns B and C (rows 1-12): header + data.
$rangeB = $ws.Range("B1:B12")
$rangeC = $ws.Range("C1:C12")

$valuesB = $rangeB.Value()
$valuesC = $rangeC.Value()

$rangeB.Value = $valuesC
$rangeC.Value = $valuesB
